$wb = $excel.ActiveWorkbook

# --- Sheet "Games" (sheet1): correct a handful of stats in row 14 ---
$games = $wb.Worksheets.Item("Games")

$games.Cells.Item(14, 5).Value = 93        # E14 Pace
$games.Cells.Item(14, 10).Value = 117.2    # J14 ORT
$games.Cells.Item(14, 13).Value = 0.616    # M14 OppeFG
$games.Cells.Item(14, 14).Value = 12.2     # N14 OppTOV
$games.Cells.Item(14, 16).Value = 0.244    # P14 OppFTR
$games.Cells.Item(14, 17).Value = 130.1    # Q14 OppORT

# --- Append the newly-played game (vs PHI on 45307) as row 43 ---
$games.Cells.Item(43, 1).Value = 42
$games.Cells.Item(43, 2).Value = 45307
$games.Cells.Item(43, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item(43, 3).Value = -1
$games.Cells.Item(43, 4).Value = 121
$games.Cells.Item(43, 5).Value = 92
$games.Cells.Item(43, 6).Value = 0.573
$games.Cells.Item(43, 7).Value = 9.1
$games.Cells.Item(43, 8).Value = 37.2
$games.Cells.Item(43, 9).Value = 0.213
$games.Cells.Item(43, 10).Value = 131.5
$games.Cells.Item(43, 11).Value = "PHI"
$games.Cells.Item(43, 12).Value = 126
$games.Cells.Item(43, 13).Value = 0.664
$games.Cells.Item(43, 14).Value = 8.1
$games.Cells.Item(43, 15).Value = 14.3
$games.Cells.Item(43, 16).Value = 0.329
$games.Cells.Item(43, 17).Value = 136.9
$games.Cells.Item(43, 18).Value = 0
$games.Cells.Item(43, 19).Value = 0

# --- Sheet "Next" (sheet2): the game vs PHI on 45307 has now been played, ---
# --- so remove it from the upcoming schedule and shift the rest up ---
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
